$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1): switch horizontal alignment from center to left.
#    (vertical="top" + bold font + border stay as-is)
# ---------------------------------------------------------------------------
$ws.Range("A1:L1").HorizontalAlignment = -4131   # xlLeft
$ws.Range("E1").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 2) Existing date columns A/B (header + body) + new row 3: left-align the
#    custom yyyy-mm-dd date format.
# ---------------------------------------------------------------------------
$ws.Range("A2:B2").HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------------
# 3) Row 2 content corrections
# ---------------------------------------------------------------------------
# C2: "08:00" text -> real time value 08:00, left aligned h:mm format
$ws.Range("C2").NumberFormat = "h:mm"
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("C2").Value = 0.33333333333333331

# D2: keep "10:00" text, but left aligned / general format
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").HorizontalAlignment = -4131
$ws.Range("D2").Value = "10:00"

# E2: category_id value 3 -> 2, left aligned integer format
$ws.Range("E2").NumberFormat = "0"
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").Value = 2

# F2: subject text, left aligned / general format (value unchanged)
$ws.Range("F2").NumberFormat = "General"
$ws.Range("F2").HorizontalAlignment = -4131
$ws.Range("F2").Value = "melakukan koordinasi dengan tim internal aplikasi 13:59"

# H2: detail_kendala text, left aligned / general format (value unchanged)
$ws.Range("H2").NumberFormat = "General"
$ws.Range("H2").HorizontalAlignment = -4131
$ws.Range("H2").Value = "Selamat sore bapak saya briptu hariawan operator program 1 polrestabes surabaya mohon petunjuk apakah ada akun untuk admin polres untuk monitoring pelaporan bhabinkamtibmas polsek jajaran pada aplikasi gugus tugas polri ketahanan pangan  ?"

# I2: respon_diberikan -> "Reset password", left aligned / general format
$ws.Range("I2").NumberFormat = "General"
$ws.Range("I2").HorizontalAlignment = -4131
$ws.Range("I2").Value = "Reset password"

# K2: priority -> "High", left aligned / general format
$ws.Range("K2").NumberFormat = "General"
$ws.Range("K2").HorizontalAlignment = -4131
$ws.Range("K2").Value = "High"

# L2: products_name -> "Photobooth", left aligned / general format
$ws.Range("L2").NumberFormat = "General"
$ws.Range("L2").HorizontalAlignment = -4131
$ws.Range("L2").Value = "Photobooth"

# ---------------------------------------------------------------------------
# 4) New columns M (status) and N (places_id)
# ---------------------------------------------------------------------------
$ws.Range("M1").Value = "status"
$ws.Range("N1").Value = "places_id"
$ws.Range("M1:N1").Font.Bold = $true
$ws.Range("M1:N1").Borders.LineStyle = 1
$ws.Range("M1:N1").HorizontalAlignment = -4131   # xlLeft
$ws.Range("M1:N1").VerticalAlignment = -4160     # xlTop

$ws.Range("N1").NumberFormat = "0"

$ws.Range("M2").NumberFormat = "General"
$ws.Range("M2").HorizontalAlignment = -4131
$ws.Range("M2").Value = "Resolved"

$ws.Range("N2").NumberFormat = "0"
$ws.Range("N2").HorizontalAlignment = -4131
$ws.Range("N2").Value = 1

$ws.Columns.Item(13).ColumnWidth = 23.71
$ws.Columns.Item(14).ColumnWidth = 20.57

# ---------------------------------------------------------------------------
# 5) New row 3 (second ticket record)
# ---------------------------------------------------------------------------
$ws.Range("A3").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").Value = 45483

$ws.Range("B3").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").Value = 45484

$ws.Range("C3").NumberFormat = "h:mm"
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("C3").Value = 0.29166666666666669

$ws.Range("D3").NumberFormat = "h:mm"
$ws.Range("D3").HorizontalAlignment = -4131
$ws.Range("D3").Value = 0.5

$ws.Range("E3").NumberFormat = "0"
$ws.Range("E3").HorizontalAlignment = -4131
$ws.Range("E3").Value = 1

$ws.Range("F3").NumberFormat = "General"
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("F3").Value = "melakukan koordinasi dengan tim internal aplikasi 13:59"

$ws.Range("H3").NumberFormat = "General"
$ws.Range("H3").HorizontalAlignment = -4131
$ws.Range("H3").Value = "Selamat sore bapak saya briptu hariawan operator program 1 polrestabes surabaya mohon petunjuk apakah ada akun untuk admin polres untuk monitoring pelaporan bhabinkamtibmas polsek jajaran pada aplikasi gugus tugas polri ketahanan pangan  ?"

$ws.Range("I3").NumberFormat = "General"
$ws.Range("I3").HorizontalAlignment = -4131
$ws.Range("I3").Value = "Reset password"

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").HorizontalAlignment = -4131
$ws.Range("J3").Value = "081319107692"

$ws.Range("K3").NumberFormat = "General"
$ws.Range("K3").HorizontalAlignment = -4131
$ws.Range("K3").Value = "Low"

$ws.Range("L3").NumberFormat = "General"
$ws.Range("L3").HorizontalAlignment = -4131
$ws.Range("L3").Value = "Gugus Pangan"

$ws.Range("M3").NumberFormat = "General"
$ws.Range("M3").HorizontalAlignment = -4131
$ws.Range("M3").Value = "New"

$ws.Range("N3").NumberFormat = "0"
$ws.Range("N3").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 6) View: select N9, scroll so column G is at the left edge
# ---------------------------------------------------------------------------
$ws.Range("N9").Select() | Out-Null
